$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition): refresh "想去人数" (want-to-go count) values ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1048
$ws1.Range("F3").Value = 380
$ws1.Range("F4").Value = 3027
$ws1.Range("F5").Value = 76
$ws1.Range("F6").Value = 630

# --- Sheet "演出" (Performance) ---
# The 2024-03-29 violin concert has passed, so the scraped feed no longer
# includes it: the 2024-03-30 "卡农" concert's details move up into row 2,
# and the old trailing row is dropped. Column A (the stale sequence index)
# is left untouched by the feed refresh.
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("B2").Value = "'2024-03-30"
$ws2.Range("C2").Value = "南宁·卡农·世界经典音乐之旅音乐会"
$ws2.Range("D2").Value = "龙堤路25号 南宁文化艺术中心"
$ws2.Range("E2").Value = "2024.03.30 20:00-03.30 21:30"
$ws2.Range("F2").Value = 19
$ws2.Range("G2").Value = 60
$ws2.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=81798"
$ws2.Range("I2").Value = "//i0.hdslb.com/bfs/openplatform/202402/Tv5lqcVn1707214065277.jpeg"
$ws2.Rows(3).Delete()

# --- Sheet "全部类型" (All types) ---
# Same 2024-03-29 row drops out; every later row's B:I content shifts up by
# one, column A stays stale, and the last (now duplicated) row is removed.
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("B2").Value = "'2024-03-30"
$ws4.Range("C2").Value = "南宁·卡农·世界经典音乐之旅音乐会"
$ws4.Range("D2").Value = "龙堤路25号 南宁文化艺术中心"
$ws4.Range("E2").Value = "2024.03.30 20:00-03.30 21:30"
$ws4.Range("F2").Value = 19
$ws4.Range("G2").Value = 60
$ws4.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=81798"
$ws4.Range("I2").Value = "//i0.hdslb.com/bfs/openplatform/202402/Tv5lqcVn1707214065277.jpeg"

$ws4.Range("B3").Value = "'2024-03-30"
$ws4.Range("C3").Value = "南宁·第一届ANE·DACG动漫嘉年华（取消）"
$ws4.Range("D3").Value = "亭洪路45号 百益上河城"
$ws4.Range("E3").Value = "2024.03.30 09:00-03.31 17:30"
$ws4.Range("F3").Value = 1048
$ws4.Range("G3").Value = "不可售"
$ws4.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=81658"
$ws4.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202403/fmbmIP421710756195423.jpeg"

$ws4.Range("B4").Value = "'2024-04-11"
$ws4.Range("C4").Value = "南宁·三月三漫次元国风动漫节"
$ws4.Range("D4").Value = "亭洪路45号 百益上河城"
$ws4.Range("E4").Value = "2024.04.11 10:00-04.12 17:00"
$ws4.Range("F4").Value = 380
$ws4.Range("G4").Value = 45
$ws4.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=83139"
$ws4.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202403/nqZxtIdd1711001896521.jpeg"

$ws4.Range("B5").Value = "'2024-05-01"
$ws4.Range("C5").Value = "南宁·2024三月三国潮动漫节（良牙春典）"
$ws4.Range("D5").Value = "民族大道106号 南宁国际会展中心"
$ws4.Range("E5").Value = "2024.05.01 09:30-05.02 17:30"
$ws4.Range("F5").Value = 3027
$ws4.Range("G5").Value = 55
$ws4.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=82416"
$ws4.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202403/b3YxmMm81711075370604.jpeg"

$ws4.Range("B6").Value = "'2024-05-19"
$ws4.Range("C6").Value = "南宁·原x穹x崩only"
$ws4.Range("D6").Value = "明秀东路157号 利泰国际大酒店"
$ws4.Range("E6").Value = "2024.05.19 10:00-05.19 17:00"
$ws4.Range("F6").Value = 76
$ws4.Range("G6").Value = 35
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=83070"
$ws4.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202403/I8tScigE1710918412731.jpeg"

$ws4.Range("B7").Value = "'2024-06-09"
$ws4.Range("C7").Value = "南宁·布谷鸟动漫展4th"
$ws4.Range("D7").Value = "亭洪路45号 百益上河城"
$ws4.Range("E7").Value = "2024.06.09 10:00-06.10 17:00"
$ws4.Range("F7").Value = 630
$ws4.Range("G7").Value = 35
$ws4.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=82241"
$ws4.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202403/uzZqZov91709281147333.jpeg"

$ws4.Rows(8).Delete()
